$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLANES_INFO")
for ($r=1; $r -le 4; $r++) {
  for ($c=1; $c -le 4; $c++) {
    $v = $ws.Cells.Item($r, $c).Value
    Write-Output "$r,$c = $v"
  }
}
